# Update the "Pais" (countries) worksheet with a newer COVID-19 data pull.
# The sheet stays sorted by column B ("Casos totales") descending, so a
# handful of neighbouring rows swap country labels as their totals change.
#
# Columns: A Pais | B Casos totales | C Nuevos casos | D Casos activos
#          E Recuperados | F Casos criticos | G Muertes hoy | H Muertes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner in A1.
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 2 de Julio de 2020 a las 10:27"

# --- Refreshed totals for countries whose ranking position is unchanged ---

# Estados Unidos (row 4)
$ws.Cells.Item(4, 2).Value = 2780152
$ws.Cells.Item(4, 3).Value = 199
$ws.Cells.Item(4, 4).Value = 1164794
$ws.Cells.Item(4, 5).Value = 1484560

# Rusia (row 6)
$ws.Cells.Item(6, 2).Value = 661165
$ws.Cells.Item(6, 3).Value = 6760
$ws.Cells.Item(6, 4).Value = 428978
$ws.Cells.Item(6, 5).Value = 222504
$ws.Cells.Item(6, 7).Value = 147
$ws.Cells.Item(6, 8).Value = 9683

# India (row 7)
$ws.Cells.Item(7, 2).Value = 605775
$ws.Cells.Item(7, 3).Value = 555
$ws.Cells.Item(7, 4).Value = 359995
$ws.Cells.Item(7, 5).Value = 227925
$ws.Cells.Item(7, 7).Value = 7
$ws.Cells.Item(7, 8).Value = 17855

# Belgica (row 30)
$ws.Cells.Item(30, 2).Value = 61598
$ws.Cells.Item(30, 3).Value = 89
$ws.Cells.Item(30, 4).Value = 17044
$ws.Cells.Item(30, 5).Value = 34793
$ws.Cells.Item(30, 7).Value = 7
$ws.Cells.Item(30, 8).Value = 9761

# Singapur (row 38)
$ws.Cells.Item(38, 2).Value = 44310
$ws.Cells.Item(38, 3).Value = 188
$ws.Cells.Item(38, 5).Value = 5273

# Polonia (row 43)
$ws.Cells.Item(43, 4).Value = 22209
$ws.Cells.Item(43, 5).Value = 11089

# Noruega (row 73)
$ws.Cells.Item(73, 2).Value = 8902
$ws.Cells.Item(73, 3).Value = 6
$ws.Cells.Item(73, 5).Value = 513

# Hungria (row 95)
$ws.Cells.Item(95, 2).Value = 4166
$ws.Cells.Item(95, 3).Value = 9
$ws.Cells.Item(95, 4).Value = 2721
$ws.Cells.Item(95, 5).Value = 858
$ws.Cells.Item(95, 7).Value = 1
$ws.Cells.Item(95, 8).Value = 587

# Estonia (row 114)
$ws.Cells.Item(114, 2).Value = 1990
$ws.Cells.Item(114, 3).Value = 1
$ws.Cells.Item(114, 5).Value = 85

# Eslovaquia (row 118)
$ws.Cells.Item(118, 2).Value = 1700
$ws.Cells.Item(118, 3).Value = 13
$ws.Cells.Item(118, 5).Value = 206

# --- Rows that swap countries because the ranking order shifted ---

# Armenia / Nigeria / Israel (rows 51-53): Armenia and Israel receive new,
# higher figures and overtake Nigeria (whose numbers are unchanged), so the
# three rows swap labels.
$ws.Cells.Item(51, 1).Value = "Armenia"
$ws.Cells.Item(51, 2).Value = 26658
$ws.Cells.Item(51, 3).Value = 593
$ws.Cells.Item(51, 4).Value = 15036
$ws.Cells.Item(51, 5).Value = 11163
$ws.Cells.Item(51, 7).Value = 6
$ws.Cells.Item(51, 8).Value = 459

$ws.Cells.Item(52, 1).Value = "Nigeria"
$ws.Cells.Item(52, 2).Value = 26484
$ws.Cells.Item(52, 4).Value = 10152
$ws.Cells.Item(52, 5).Value = 15729
$ws.Cells.Item(52, 8).Value = 603

$ws.Cells.Item(53, 1).Value = "Israel"
$ws.Cells.Item(53, 2).Value = 26452
$ws.Cells.Item(53, 3).Value = 195
$ws.Cells.Item(53, 4).Value = 17481
$ws.Cells.Item(53, 5).Value = 8647
$ws.Cells.Item(53, 7).Value = 2
$ws.Cells.Item(53, 8).Value = 324

# Malaui overtakes Cabo Verde (rows 125-126).
$ws.Cells.Item(125, 1).Value = "Malaui"
$ws.Cells.Item(125, 2).Value = 1342
$ws.Cells.Item(125, 3).Value = 77
$ws.Cells.Item(125, 4).Value = 271
$ws.Cells.Item(125, 5).Value = 1055
$ws.Cells.Item(125, 8).Value = 16

$ws.Cells.Item(126, 1).Value = "Cabo Verde"
$ws.Cells.Item(126, 2).Value = 1267
$ws.Cells.Item(126, 4).Value = 629
$ws.Cells.Item(126, 5).Value = 623
$ws.Cells.Item(126, 8).Value = 15

# Birmania overtakes Comoras (rows 161-162).
$ws.Cells.Item(161, 1).Value = "Birmania"
$ws.Cells.Item(161, 2).Value = 304
$ws.Cells.Item(161, 3).Value = 1
$ws.Cells.Item(161, 4).Value = 222
$ws.Cells.Item(161, 5).Value = 76
$ws.Cells.Item(161, 8).Value = 6

$ws.Cells.Item(162, 1).Value = "Comoras"
$ws.Cells.Item(162, 4).Value = 200
$ws.Cells.Item(162, 5).Value = 96
$ws.Cells.Item(162, 8).Value = 7

# Laos and Santa Lucia are tied at 19 cases; Laos moves ahead (rows 203-204).
$ws.Cells.Item(203, 1).Value = "Laos"
$ws.Cells.Item(204, 1).Value = "Santa Lucia"
